# Atualização dos endereços no import-2.xlsx
#
# Updates the ENDERECOS (address) column (K) for the 5 "JURIDICA" (company)
# rows that previously all shared the placeholder "...;75100100;Distrito
# Industrial;Anápolis;GO" address, replacing them with real São Paulo
# addresses. Also restores the active-cell selection left behind by the
# edit (K1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Silva & Cia Ltda.
$ws.Range("K3").Value = "Avenida do estado;1000;01025020;Centro;São Paulo;SP"

# Row 4 - Costa Transportes
$ws.Range("K4").Value = "Avenida Antártica;s/n;05003020;Água Branca;São Paulo;SP"

# Row 5 - Dias Comércio de Alimentos
$ws.Range("K5").Value = "Estrada do carrapicho;4;05275015;Anhanguera;São Paulo;SP"

# Row 6 - Moreira & Filhos ME
$ws.Range("K6").Value = "Rua Ágatha Cristie;10;04875160;Chácara Santo Amaro;São Paulo;SP"

# Row 7 - F&E Consultoria (street/number unchanged, only CEP/bairro/cidade/UF updated)
$ws.Range("K7").Value = "Rua janeiros;1000;04116000;Jardim Vila Mariana;São Paulo;SP"

# The saved workbook's selection moved from A1 to K1 (last touched column).
$null = $ws.Range("K1").Select()
